$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Force a value to be stored as text even if it looks like a number or
    # date (e.g. "54345" or "2022-11-17"), then drop back to the Normal
    # style so no stray number-format override is left on the cell.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# ----- Sheet 1: "My Foods" -----
$ws1 = $wb.Worksheets.Item(1)

$foodRows = @(
    @("dsasdds", 91.033, 25, 33, 33, 0, 33, 0, 529),
    @("54345",   91.033, 25, 33, 33, 0, 33, 0, 529),
    @("dfdsf",   91.033, 25, 33, 33, 0, 33, 0, 529),
    @("987gg",   91.033, 25, 33, 33, 0, 33, 0, 529),
    @("aaaa",    91.033, 25, 33, 33, 0, 33, 0, 529),
    @("bbbb",    91.033, 25, 33, 33, 0, 33, 0, 529),
    @("ccc",     91.033, 25, 33, 33, 0, 33, 0, 529),
    @("ddd",     91.033, 25, 33, 33, 0, 33, 0, 529),
    @("eeeee",   91.033, 25, 33, 33, 0, 33, 0, 529),
    @("fffs",    91.033, 25, 33, 33, 0, 33, 0, 529)
)

$r = 5
foreach ($row in $foodRows) {
    Set-TextValue $ws1.Cells.Item($r, 1) $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $ws1.Cells.Item($r, 6).Value = $row[5]
    $ws1.Cells.Item($r, 7).Value = $row[6]
    $ws1.Cells.Item($r, 8).Value = $row[7]
    $ws1.Cells.Item($r, 9).Value = $row[8]
    $r = $r + 1
}

# ----- Sheet 2: "My Meal Entries" -----
$ws2 = $wb.Worksheets.Item(2)

$mealRows = @(
    @("2022-11-17", "Sherbet", 98.760046, 1.1, 1.16),
    @("2022-11-17", "Sherbet", 98.760046, 1.1, 1.16),
    @("2022-11-17", "987gg", 91.033, 25, 33)
)

$r = 7
foreach ($row in $mealRows) {
    Set-TextValue $ws2.Cells.Item($r, 1) $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
